# Replace placeholder serial numbers with real/example values and fix the
# swapped "Backend" values on the "Connections" sheet (include usage example
# in API intro).
#
# Writes are ordered so that new shared-string entries are introduced in the
# same sequence as the canonical edit: D7D59860, 8BB438A7, 51957DB9,
# 2CAD4BC6, USB::0x7B1A::0x0381::D7D59860, MY24339283,
# USB::0x9D0A::0x2217::MY24339283.

$wb = $excel.ActiveWorkbook

$eq = $wb.Worksheets.Item("Equipment")
$conn = $wb.Worksheets.Item("Connections")

$eq.Range("C2").Value = "D7D59860"
$eq.Range("C4").Value = "8BB438A7"
$eq.Range("C5").Value = "51957DB9"
$eq.Range("C6").Value = "2CAD4BC6"
$conn.Range("E2").Value = "USB::0x7B1A::0x0381::D7D59860"
$eq.Range("C7").Value = "MY24339283"
$conn.Range("E7").Value = "USB::0x9D0A::0x2217::MY24339283"

$conn.Range("C2").Value = "D7D59860"
$conn.Range("D3").Value = "MSL"
$conn.Range("C4").Value = "8BB438A7"
$conn.Range("D4").Value = "PyVISA"
$conn.Range("C5").Value = "51957DB9"
$conn.Range("C6").Value = "2CAD4BC6"
$conn.Range("C7").Value = "MY24339283"
